$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Alternative")
$ws.Range("J2").Value = 9.789999999999999
$ws.Range("K2").Value = 9.789999999999999
$ws.Range("L2").Value = 12.28
$ws.Range("M2").Value = 4.83
$ws.Range("N2").Value = -10.22
$ws.Range("O2").Value = 0.37
$ws.Range("J3").Value = 1.53
$ws.Range("K3").Value = 9.359999999999999
$ws.Range("L3").Value = 8.380000000000001
$ws.Range("M3").Value = 7.06
$ws.Range("N3").Value = -4.06
$ws.Range("O3").Value = 0.49
$ws.Range("J4").Value = 1.26
$ws.Range("K4").Value = 3.18
$ws.Range("L4").Value = 7.44
$ws.Range("N4").Value = -9.039999999999999
$ws.Range("O4").Value = -0.27
$ws.Range("J5").Value = -8.550000000000001
$ws.Range("K5").Value = -8.550000000000001
$ws.Range("L5").Value = 44.19
$ws.Range("M5").Value = 26.64
$ws.Range("N5").Value = -82.28
$ws.Range("O5").Value = -0.31
$ws.Range("J6").Value = -0.25
$ws.Range("K6").Value = 2.13
$ws.Range("L6").Value = 6.54
$ws.Range("M6").Value = 7.88
$ws.Range("N6").Value = -8.550000000000001
$ws.Range("O6").Value = -0.47
$ws.Range("J7").Value = 1.52
$ws.Range("K7").Value = 1.52
$ws.Range("L7").Value = 5.42
$ws.Range("M7").Value = 5.85
$ws.Range("N7").Value = -7.34
$ws.Range("O7").Value = -0.68

$ws = $wb.Worksheets.Item("Bond")
$ws.Range("J2").Value = -3.74
$ws.Range("K2").Value = 11.5
$ws.Range("L2").Value = 10.81
$ws.Range("M2").Value = 12.83
$ws.Range("N2").Value = -6.38
$ws.Range("O2").Value = 0.58
$ws.Range("N3").Value = 3.7
$ws.Range("O3").Value = -2.84
$ws.Range("J4").Value = -2.15
$ws.Range("K4").Value = 3.28
$ws.Range("L4").Value = 17.89
$ws.Range("M4").Value = 15.8
$ws.Range("N4").Value = -26.26
$ws.Range("J5").Value = -3.63
$ws.Range("K5").Value = 2.38
$ws.Range("L5").Value = 11.56
$ws.Range("M5").Value = 13.11
$ws.Range("N5").Value = -16.61

$ws = $wb.Worksheets.Item("Equity")
$ws.Range("J2").Value = 12.92
$ws.Range("K2").Value = 14.56
$ws.Range("L2").Value = 19.75
$ws.Range("M2").Value = 11.12
$ws.Range("N2").Value = -17.39
$ws.Range("J3").Value = 8.779999999999999
$ws.Range("K3").Value = 11.98
$ws.Range("L3").Value = 16.71
$ws.Range("M3").Value = 8.640000000000001
$ws.Range("N3").Value = -15.37
$ws.Range("J4").Value = 4.22
$ws.Range("K4").Value = 6.48
$ws.Range("L4").Value = 14.83
$ws.Range("M4").Value = 10.84
$ws.Range("N4").Value = -17.65
$ws.Range("O4").Value = 0.09
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 8.74
$ws.Range("L5").Value = 10.64
$ws.Range("M5").Value = 7.52
$ws.Range("N5").Value = -8.779999999999999
$ws.Range("J6").Value = 1.3
$ws.Range("K6").Value = 5.15
$ws.Range("L6").Value = 26.45
$ws.Range("M6").Value = 17.11
$ws.Range("N6").Value = -39.36
$ws.Range("O6").Value = (0.0 * -1)
$ws.Range("J7").Value = 5.31
$ws.Range("K7").Value = 8.67
$ws.Range("L7").Value = 15.35
$ws.Range("M7").Value = 10.65
$ws.Range("N7").Value = -16.25
$ws.Range("O7").Value = 0.23
$ws.Range("J8").Value = 1.58
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 14.48
$ws.Range("M8").Value = 11.9
$ws.Range("N8").Value = -19.05
$ws.Range("L9").Value = 13.66
$ws.Range("M9").Value = 16.13
$ws.Range("N9").Value = -23.6
$ws.Range("O9").Value = -0.45
$ws.Range("J10").Value = 2.35
$ws.Range("K10").Value = 7.4
$ws.Range("L10").Value = 26.07
$ws.Range("M10").Value = 16.52
$ws.Range("N10").Value = -35.76
$ws.Range("O10").Value = 0.08
$ws.Range("J11").Value = -1.11
$ws.Range("K11").Value = 11.83
$ws.Range("L11").Value = 7.96
$ws.Range("M11").Value = 9.119999999999999
$ws.Range("N11").Value = -1.36
$ws.Range("O11").Value = 0.83

